# Regenerate column G ("K") values in the save_data sheet.
# The commit replaces the previous "Strike#"-derived counts that were
# stored in column G with newly (re)computed K values, alongside a
# (non-visible-in-this-sheet) regeneration of std/mean and s_vals.
# Here we simply write the new, already-computed K values into column G
# for each data row (rows 2-40), leaving every other cell untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 1
$ws.Cells.Item(3, 7).Value = 1
$ws.Cells.Item(4, 7).Value = 3
$ws.Cells.Item(5, 7).Value = 2
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(7, 7).Value = 2
$ws.Cells.Item(8, 7).Value = 2
$ws.Cells.Item(9, 7).Value = 2
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(12, 7).Value = 1
$ws.Cells.Item(13, 7).Value = 2
$ws.Cells.Item(14, 7).Value = 1
$ws.Cells.Item(15, 7).Value = 2
$ws.Cells.Item(16, 7).Value = 3
$ws.Cells.Item(17, 7).Value = 2
$ws.Cells.Item(18, 7).Value = 2
$ws.Cells.Item(19, 7).Value = 1
$ws.Cells.Item(20, 7).Value = 2
$ws.Cells.Item(21, 7).Value = 1
$ws.Cells.Item(22, 7).Value = 3
$ws.Cells.Item(23, 7).Value = 1
$ws.Cells.Item(24, 7).Value = 2
$ws.Cells.Item(25, 7).Value = 1
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(27, 7).Value = 1
$ws.Cells.Item(28, 7).Value = 3
$ws.Cells.Item(29, 7).Value = 0
$ws.Cells.Item(30, 7).Value = 1
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(32, 7).Value = 4
$ws.Cells.Item(33, 7).Value = 1
$ws.Cells.Item(34, 7).Value = 1
$ws.Cells.Item(35, 7).Value = 0
$ws.Cells.Item(36, 7).Value = 2
$ws.Cells.Item(37, 7).Value = 2
$ws.Cells.Item(38, 7).Value = 0
$ws.Cells.Item(39, 7).Value = 2
$ws.Cells.Item(40, 7).Value = 0
